# Swap the data content of rows 30 and 31 (columns C:I) on the active sheet,
# then move the selection to E37 (matching the author's saved selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for row 30 and row 31, columns C..I
$cols = @("C", "D", "E", "F", "G", "H", "I")

$row30 = @{}
$row31 = @{}
foreach ($col in $cols) {
    $row30[$col] = $ws.Range("$col" + "30").Value2
    $row31[$col] = $ws.Range("$col" + "31").Value2
}

# Write row 31's original content into row 30, and row 30's original content into row 31
foreach ($col in $cols) {
    $ws.Range("$col" + "30").Value2 = $row31[$col]
    $ws.Range("$col" + "31").Value2 = $row30[$col]
}

# Update the stored selection to match the author's saved view state
$ws.Range("E37").Select()
